$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Runs")

# Insert a new row at position 50 (shifts "Run058".."Run8" down by one row,
# i.e. old row 50 becomes row 51, ..., old row 64 becomes row 65).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the "Run058_final" results.
$ws.Range("A50").Value = 0
$ws.Range("B50").Value = "Run058_final"
$ws.Range("C50").Value = "0.852 (0.261)"
$ws.Range("D50").Value = "0.654 (0.195)"
$ws.Range("E50").Value = "0.941 (0.237)"
$ws.Range("F50").Value = "0.918 (0.182)"
$ws.Range("G50").Value = "0.799 (0.156)"
$ws.Range("H50").Value = "0.972 (0.166)"
$ws.Range("I50").Value = "0.934 (0.180)"
$ws.Range("J50").Value = "0.767 (0.197)"
$ws.Range("K50").Value = "0.978 (0.147)"
$ws.Range("L50").Value = "0.645 (0.383)"
$ws.Range("M50").Value = "0.645 (0.238)"
$ws.Range("N50").Value = "0.646 (0.480)"
$ws.Range("O50").Value = "0.821 (0.331)"
$ws.Range("P50").Value = "0.687 (0.322)"
$ws.Range("Q50").Value = "0.889 (0.315)"
$ws.Range("R50").Value = "0.852 (0.310)"
$ws.Range("S50").Value = "0.703 (0.209)"
$ws.Range("T50").Value = "0.885 (0.320)"
$ws.Range("U50").Value = "0.818 (0.330)"
$ws.Range("V50").Value = "0.624 (0.278)"
$ws.Range("W50").Value = "0.884 (0.321)"
$ws.Range("X50").Value = "0.929 (0.215)"
$ws.Range("Y50").Value = "0.744 (0.229)"
$ws.Range("Z50").Value = "0.960 (0.197)"
$ws.Range("AA50").Value = "0.755 (0.415)"
$ws.Range("AB50").Value = "0.588 (0.204)"
$ws.Range("AC50").Value = "0.766 (0.424)"
$ws.Range("AD50").Value = "0.799 (0.325)"
$ws.Range("AE50").Value = "0.641 (0.238)"
$ws.Range("AF50").Value = "0.873 (0.334)"
$ws.Range("AG50").Value = "0.905 (0.233)"
$ws.Range("AH50").Value = "0.756 (0.231)"
$ws.Range("AI50").Value = "0.952 (0.213)"
$ws.Range("AJ50").Value = "0.841 (0.333)"
$ws.Range("AK50").Value = "0.716 (0.212)"
$ws.Range("AL50").Value = "0.862 (0.345)"
$ws.Range("AM50").Value = "0.742 (0.326)"
$ws.Range("AN50").Value = "0.702 (0.218)"
$ws.Range("AO50").Value = "0.786 (0.410)"

# Row-insert re-materializes the originally-blank trailing cells of the old
# "Run063" row (now shifted to row 56, columns O:AO) as empty strings; clear
# them back to true blanks to match the source data (which never had values there).
$ws.Range("O56:AO56").ClearContents()
